# Update Commodity, Process and Process-Commodity sheets with new dev values
# (haag_wgs84), as per commit "updated runmoosh.py and data.xlsx with devs
# from haag_wgs84".

$wb = $excel.ActiveWorkbook

$wsCommodity = $wb.Worksheets.Item("Commodity")
$wsProcess = $wb.Worksheets.Item("Process")
$wsProcessCommodity = $wb.Worksheets.Item("Process-Commodity")

# --- Process-Commodity sheet ---
$wsProcessCommodity.Range("D13").Value = 0.99

# --- Process sheet ---
$wsProcess.Range("C2").Value = 0.2

$wsProcess.Range("C4").Value = 0.2
$wsProcess.Range("F4").Value = 1000

$wsProcess.Range("F7").Value = 999990

$wsProcess.Range("C8").Value = 0.3
$wsProcess.Range("F8").Value = 1000

$wsProcess.Select()
$wsProcess.Range("D3:D4").Select()

# --- Commodity sheet ---
$wsCommodity.Range("F2").Value = 0.001
$wsCommodity.Range("H2").Value = 750000

$wsCommodity.Range("F3").Value = 0.0001
$wsCommodity.Range("H3").Value = 500000

$wsCommodity.Range("C4").Value = 1500
$wsCommodity.Range("D4").Value = 0.01
$wsCommodity.Range("F4").Value = 1

$wsCommodity.Range("C5").Value = 1500
$wsCommodity.Range("D5").Value = 0.01
$wsCommodity.Range("F5").Value = 2

# Commodity must remain the active sheet/tab (it was the only one with
# tabSelected originally), so select it - and its target cell - last.
$wsCommodity.Select()
$wsCommodity.Range("C4").Select()
